# Append the new YIMBY news row (row 2) to Sheet1, which currently only
# has the header row (row 1: Title, URL, Summary, Published, raw_date,
# Source, Feed, content_preview).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "45–40 Vernon Boulevard Tops Out in Long Island City, Queens"
$ws.Range("B2").Value = "https://newyorkyimby.com/2025/11/45-40-vernon-boulevard-tops-out-in-long-island-city-queens.html"
$ws.Range("C2").Value = 'Construction has topped out on 45–40 Vernon Boulevard, a 23-story residential building in <a href="https://newyorkyimby.com/neighborhoods/long-island-city">Long Island City</a>, Queens. Designed by Archimaera and developed and built by ZD Jasper, the 262-foot-tall structure will span 192,500 square feet and yield 226 units along with lower-level commercial space. The property is located near the intersection of Vernon Boulevard and 46th Avenue.'
$ws.Range("D2").Value = "2025-11-25T12:30:14+00:00"
$ws.Range("E2").Value = "Tue, 25 Nov 2025 12:30:14 +0000"
$ws.Range("F2").Value = "YIMBY"
$ws.Range("G2").Value = "YIMBY - Long Island City"

# content_preview is blank for this item (source feed had no preview
# text). An empty assignment keeps column H in the used range without
# introducing any visible content.
$ws.Range("H2").Value = ""
